$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.503.76"
$ws.Range("E2").Value = "  +3.81%  "

$ws.Range("D3").Value = "1.839.54"
$ws.Range("E3").Value = "  +2.77%  "

$ws.Range("E4").Value = "  +2.64%  "

$ws.Range("D5").Value = "318.99"
$ws.Range("E5").Value = "  +3.77%  "

$ws.Range("D6").Value = "1.025"
$ws.Range("E6").Value = "  +2.34%  "

$ws.Range("D7").Value = "0.4371"
$ws.Range("E7").Value = "  +2.86%  "

$ws.Range("D8").Value = "0.3733"
$ws.Range("E8").Value = "  +3.32%  "

$ws.Range("D9").Value = "0.07377"
$ws.Range("E9").Value = "  +2.98%  "

$ws.Range("E10").Value = "  +3.04%  "

$ws.Range("D11").Value = "21.50"
$ws.Range("E11").Value = "  +4.56%  "

$ws.Range("D12").Value = "1.844.05"
$ws.Range("E12").Value = "  +4.50%  "

$ws.Range("D13").Value = "5.495"
$ws.Range("E13").Value = "  +4.43%  "

$ws.Range("D14").Value = "6.679"
$ws.Range("E14").Value = "  +2.98%  "

$ws.Range("D15").Value = "0.07148"
$ws.Range("E15").Value = "  +3.48%  "

$ws.Range("D16").Value = "82.62"
$ws.Range("E16").Value = "  +4.16%  "

$ws.Range("D17").Value = "1.031"
$ws.Range("E17").Value = "  +2.43%  "

$ws.Range("D18").Value = "0.000009001"
$ws.Range("E18").Value = "  +2.78%  "

$ws.Range("D19").Value = "1.026"
$ws.Range("E19").Value = "  +2.16%  "

$ws.Range("D20").Value = "15.40"
$ws.Range("E20").Value = "  +2.99%  "

$ws.Range("D21").Value = "27.513.71"
$ws.Range("E21").Value = "  +3.87%  "

$ws.Range("D22").Value = "5.253"
$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("D23").Value = "11.18"
$ws.Range("E23").Value = "  +1.19%  "

$ws.Range("D24").Value = "2.058.60"
$ws.Range("E24").Value = "  +3.86%  "

$ws.Range("D25").Value = "157.21"
$ws.Range("E25").Value = "  +3.15%  "

$ws.Range("D26").Value = "1.928"
$ws.Range("E26").Value = "  +6.62%  "

$ws.Range("E27").Value = "  +3.19%  "

$ws.Range("D28").Value = "5.245"
$ws.Range("E28").Value = "  +2.48%  "

$ws.Range("E29").Value = "  +3.85%  "

$ws.Range("D30").Value = "115.87"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("D31").Value = "0.09097"

$ws.Range("D32").Value = "1.205"
$ws.Range("E32").Value = "  +6.00%  "

$ws.Range("D33").Value = "0.7678"
$ws.Range("E33").Value = "  +4.23%  "

$ws.Range("D34").Value = "4.498"
$ws.Range("E34").Value = "  +3.57%  "

$ws.Range("D35").Value = "2.869"
$ws.Range("E35").Value = "  +4.49%  "

$ws.Range("D36").Value = "1.027"
$ws.Range("E36").Value = "  +2.69%  "

$ws.Range("D37").Value = "1.140"
$ws.Range("E37").Value = "  +2.57%  "

$ws.Range("D38").Value = "0.01971"
$ws.Range("E38").Value = "  +4.02%  "

$ws.Range("E39").Value = "  +1.74%  "

$ws.Range("E40").Value = "  +4.06%  "

$ws.Range("D41").Value = "2.780"
$ws.Range("E41").Value = "  +6.37%  "

$ws.Range("D42").Value = "0.1669"
$ws.Range("E42").Value = "  +3.11%  "

$ws.Range("D43").Value = "6.652"
$ws.Range("E43").Value = "  +4.31%  "

$ws.Range("D44").Value = "8.547"
$ws.Range("E44").Value = "  +4.35%  "

$ws.Range("D45").Value = "108.81"
$ws.Range("E45").Value = "  +3.40%  "

$ws.Range("D46").Value = "10.53"
$ws.Range("E46").Value = "  +2.70%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "1.712"
$ws.Range("E47").Value = "  +4.74%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.4643"
$ws.Range("E48").Value = "  +3.03%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.900"
$ws.Range("E49").Value = "  +7.90%  "

$ws.Range("D50").Value = "0.06343"
$ws.Range("E50").Value = "  +2.27%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "39.50"
$ws.Range("E51").Value = "  +7.08%  "
